$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the weighting factors in column E for rows 4 and 6
$ws.Range("E4").Value = 0.5
$ws.Range("E6").Value = 0.5

# Update the active cell selection to match the edited workbook
$ws.Range("F17").Select()
